# Generate Report for Handoff
#
# The two tracked files in this localization-status report are
#   b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md
#   fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md
#
# This edit re-orders the rows on every sheet (fd27d0a9 now sorts first,
# b64c069e second) and refreshes b64c069e's handoff status: it moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with new
# handoff timestamps, a refreshed target-xlf name and a populated
# "Error Detail" explaining the handback file is stale.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md"
$ov.Range("B2").Value = "e2e\fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md"

$ov.Range("A3").Value = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md"
$ov.Range("B3").Value = "e2e\b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-18 00:45:45"

# Hyperlinks on B2/B3 keep pointing at the same targets (rId2 -> b64c069e,
# rId3 -> fd27d0a9) but the two display texts now swap along with the data.
# (TextToDisplay can't be mutated in place on this host -- it appends a new
# hyperlink record instead of replacing -- so delete + recreate instead.)
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40b961c28f3fd911bb36d9cfcf6be053a1ef985e/e2e/b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md", $null, $null, "e2e\fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40b961c28f3fd911bb36d9cfcf6be053a1ef985e/e2e/fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md", $null, $null, "e2e\b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md"
$zh.Range("G2").Value = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.931a24320b095df14becc982733cbba216afdc2e.zh-cn.xlf"
$zh.Range("I2").Value = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md"
$zh.Range("J2").Value = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.931a24320b095df14becc982733cbba216afdc2e.zh-cn.xlf"

$zh.Range("A3").Value = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.f5bd270e4a6c3bf1e09824f3b16df25fac5d5685.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-18 00:45:40"
$zh.Range("I3").Value = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md"
$zh.Range("J3").Value = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.f5bd270e4a6c3bf1e09824f3b16df25fac5d5685.zh-cn.xlf"
$zh.Range("P3").Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40b961c28f3fd911bb36d9cfcf6be053a1ef985e/e2e/b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/78c48669a7db30259decd2748e9dc98fd0fc08d3/e2e/b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md.'

$zh.Hyperlinks.Item(1).TextToDisplay = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md"
$zh.Hyperlinks.Item(2).TextToDisplay = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md"
$zh.Hyperlinks.Item(3).TextToDisplay = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md"
$zh.Hyperlinks.Item(4).TextToDisplay = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md"

$zh.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md"
$de.Range("G2").Value = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.931a24320b095df14becc982733cbba216afdc2e.de-de.xlf"
$de.Range("I2").Value = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md"
$de.Range("J2").Value = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.931a24320b095df14becc982733cbba216afdc2e.de-de.xlf"

$de.Range("A3").Value = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.f5bd270e4a6c3bf1e09824f3b16df25fac5d5685.de-de.xlf"
$de.Range("H3").Value = "2016-08-18 00:45:45"
$de.Range("I3").Value = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md"
$de.Range("J3").Value = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.f5bd270e4a6c3bf1e09824f3b16df25fac5d5685.de-de.xlf"
$de.Range("P3").Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40b961c28f3fd911bb36d9cfcf6be053a1ef985e/e2e/b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/78c48669a7db30259decd2748e9dc98fd0fc08d3/e2e/b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md.'

$de.Hyperlinks.Item(1).TextToDisplay = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md"
$de.Hyperlinks.Item(2).TextToDisplay = "fd27d0a9-bd0c-4e03-8fe1-d97f19ea9a19.md"
$de.Hyperlinks.Item(3).TextToDisplay = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md"
$de.Hyperlinks.Item(4).TextToDisplay = "b64c069e-93b1-4e0a-8254-3dcf1aa542bb.md"

$de.Columns.Item(16).ColumnWidth = 39.17
